$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the "Date" value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-07-25T07:22:51+00:00"

# --- Elements sheet: update canonical terminology URLs ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("Z3").Value = "https://mos.esante.gouv.fr/NOS/TRE_R14-TypeDiplome/FHIR/TRE-R14-TypeDiplome?vs"
$elements.Range("Z4").Value = "https://mos.esante.gouv.fr/NOS/TRE_R16-LieuFormation/FHIR/TRE-R16-LieuFormation?vs"
$elements.Range("Z7").Value = "https://mos.esante.gouv.fr/NOS/TRE_R56-Attestation/FHIR/TRE-R56-Attestation?vs"

# --- Widen column Z to fit the new (longer) URL text ---
# (ColumnWidth is quantized to 1/6-character steps by this host; 70 is the
# closest input that reproduces the saved OOXML width of ~70.8555 used by
# the original authoring tool.)
$elements.Columns.Item(26).ColumnWidth = 70
